$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2565098.5
$ws.Range("J17").Value = 2565098.5
$ws.Range("L17").Value = 7695295.5
$ws.Range("N17").Value = -7695631.5

$ws.Range("H28").Value = 866.1539
$ws.Range("I28").Value = 1300
$ws.Range("J28").Value = 673.3333
$ws.Range("K28").Value = 1300
$ws.Range("L28").Value = 673.3333
$ws.Range("M28").Value = -815
$ws.Range("N28").Value = -1643.3333

$ws.Range("H33").Value = 18306.666
$ws.Range("I33").Value = 19206.941
$ws.Range("K33").Value = 19206.941
$ws.Range("M33").Value = -18977.941

$ws.Range("H64").Value = 9402.588
$ws.Range("I64").Value = 3263.5715
$ws.Range("K64").Value = 3263.5715
$ws.Range("M64").Value = -3015.5715

$ws.Range("H67").Value = 9402.588
$ws.Range("I67").Value = 3263.5715
$ws.Range("K67").Value = 3263.5715
$ws.Range("M67").Value = -2405.5715

$ws.Range("H70").Value = 4175.125
$ws.Range("I70").Value = 1733.6666
$ws.Range("J70").Value = 4523.905
$ws.Range("K70").Value = 5200.9998
$ws.Range("L70").Value = 13571.715
$ws.Range("M70").Value = -4930.9998
$ws.Range("N70").Value = -14111.715

$ws.Range("H73").Value = 4175.125
$ws.Range("I73").Value = 1733.6666
$ws.Range("J73").Value = 4523.905
$ws.Range("K73").Value = 5200.9998
$ws.Range("L73").Value = 13571.715
$ws.Range("M73").Value = -4264.9998
$ws.Range("N73").Value = -15443.715

$ws.Range("H96").Value = 354.53845
$ws.Range("I96").Value = 403.36365
$ws.Range("K96").Value = 1210.09095
$ws.Range("M96").Value = 162.90905

$ws.Range("H106").Value = 1847.4
$ws.Range("I106").Value = 1759.25
$ws.Range("K106").Value = 1759.25
$ws.Range("M106").Value = -1128.25

$ws.Range("H116").Value = 4017.0908
$ws.Range("J116").Value = 4463.6665
$ws.Range("L116").Value = 4463.6665
$ws.Range("N116").Value = -11347.6665

$ws.Range("H132").Value = 4004.3784
$ws.Range("I132").Value = 873.14813
$ws.Range("K132").Value = 2619.44439
$ws.Range("M132").Value = -89.44439000000011

$ws.Range("H141").Value = 2674.9333
$ws.Range("I141").Value = 2476.08
$ws.Range("K141").Value = 7428.24
$ws.Range("M141").Value = -2248.24

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 555
$ws.Range("I16").Value = 555
$ws.Range("K16").Value = 555
$ws.Range("M16").Value = -268

$ws.Range("H74").Value = 3955.8965
$ws.Range("I74").Value = 3522.1765
$ws.Range("J74").Value = 4570.3335
$ws.Range("K74").Value = 3522.1765
$ws.Range("L74").Value = 4570.3335
$ws.Range("M74").Value = -2648.1765
$ws.Range("N74").Value = -6318.3335

$ws.Range("H77").Value = 3955.8965
$ws.Range("I77").Value = 3522.1765
$ws.Range("J77").Value = 4570.3335
$ws.Range("K77").Value = 17610.8825
$ws.Range("L77").Value = 22851.6675
$ws.Range("M77").Value = -13242.8825
$ws.Range("N77").Value = -31587.6675

$ws.Range("H102").Value = 5579.7646
$ws.Range("I102").Value = 4204
$ws.Range("J102").Value = 12000
$ws.Range("K102").Value = 4204
$ws.Range("L102").Value = 12000
$ws.Range("M102").Value = -2582
$ws.Range("N102").Value = -15244

$ws.Range("H132").Value = 5311.41
$ws.Range("I132").Value = 4489.8857
$ws.Range("K132").Value = 13469.6571
$ws.Range("M132").Value = -10939.6571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5222.1665
$ws.Range("I16").Value = 3833.25
$ws.Range("K16").Value = 3833.25
$ws.Range("M16").Value = -3546.25

$ws.Range("H31").Value = 5669.091
$ws.Range("I31").Value = 3627.7856
$ws.Range("J31").Value = 7173.2104
$ws.Range("K31").Value = 3627.7856
$ws.Range("L31").Value = 7173.2104
$ws.Range("M31").Value = -3332.7856
$ws.Range("N31").Value = -7763.2104

$ws.Range("H34").Value = 5669.091
$ws.Range("I34").Value = 3627.7856
$ws.Range("J34").Value = 7173.2104
$ws.Range("K34").Value = 3627.7856
$ws.Range("L34").Value = 7173.2104
$ws.Range("M34").Value = -3425.7856
$ws.Range("N34").Value = -7577.2104

$ws.Range("H86").Value = 10076.6
$ws.Range("I86").Value = 10095.75
$ws.Range("K86").Value = 10095.75
$ws.Range("M86").Value = -8972.75

$ws.Range("H89").Value = 10076.6
$ws.Range("I89").Value = 10095.75
$ws.Range("K89").Value = 50478.75
$ws.Range("M89").Value = -44862.75

$ws.Range("H113").Value = 5222.1665
$ws.Range("I113").Value = 3833.25
$ws.Range("K113").Value = 3833.25
$ws.Range("M113").Value = -1663.25

$ws.Range("H135").Value = 81390.78
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 81390.78
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 81390.78
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -91530.78

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 909765.2
$ws.Range("I117").Value = 484.5
$ws.Range("J117").Value = 1111827.5
$ws.Range("K117").Value = 1453.5
$ws.Range("L117").Value = 3335482.5
$ws.Range("M117").Value = 1988.5
$ws.Range("N117").Value = -3342366.5

$ws.Range("H121").Value = 7692877
$ws.Range("J121").Value = 16667379
$ws.Range("L121").Value = 50002137
$ws.Range("N121").Value = -50004757

$ws.Range("H129").Value = 6660.5884
$ws.Range("I129").Value = 1646.375
$ws.Range("J129").Value = 11117.667
$ws.Range("K129").Value = 4939.125
$ws.Range("L129").Value = 33353.001
$ws.Range("M129").Value = 60.875
$ws.Range("N129").Value = -43353.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2457
$ws.Range("I102").Value = 1924.4445
$ws.Range("K102").Value = 1924.4445
$ws.Range("M102").Value = -302.4445000000001

$ws.Range("H122").Value = 3904.0908
$ws.Range("I122").Value = 3070.7144
$ws.Range("J122").Value = 5362.5
$ws.Range("K122").Value = 9212.143199999999
$ws.Range("L122").Value = 16087.5
$ws.Range("M122").Value = -6762.143199999999
$ws.Range("N122").Value = -20987.5

$ws.Range("H132").Value = 7099.4
$ws.Range("I132").Value = 6549.9287
$ws.Range("K132").Value = 19649.7861
$ws.Range("M132").Value = -17119.7861

$ws.Range("H135").Value = 70000
$ws.Range("J135").Value = 70000
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140

$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1734.8462
$ws.Range("J16").Value = 1835.8
$ws.Range("L16").Value = 1835.8
$ws.Range("N16").Value = -2175.8

$ws.Range("H136").Value = 8061.6
$ws.Range("I136").Value = 7160.4287
$ws.Range("K136").Value = 21481.2861
$ws.Range("M136").Value = -18931.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 8031.1113
$ws.Range("I25").Value = 8030
$ws.Range("J25").Value = 8035
$ws.Range("K25").Value = 8030
$ws.Range("L25").Value = 8035
$ws.Range("M25").Value = -7737
$ws.Range("N25").Value = -8621

$ws.Range("H70").Value = 37398
$ws.Range("I70").Value = 25000
$ws.Range("J70").Value = 40497.5
$ws.Range("K70").Value = 25000
$ws.Range("L70").Value = 40497.5
$ws.Range("M70").Value = -24685
$ws.Range("N70").Value = -41127.5

$ws.Range("H73").Value = 37398
$ws.Range("I73").Value = 25000
$ws.Range("J73").Value = 40497.5
$ws.Range("K73").Value = 25000
$ws.Range("L73").Value = 40497.5
$ws.Range("M73").Value = -23908
$ws.Range("N73").Value = -42681.5

$ws.Range("H96").Value = 1459.5
$ws.Range("I96").Value = 1550
$ws.Range("K96").Value = 1550
$ws.Range("M96").Value = -177

$ws.Range("H100").Value = 1168.8182
$ws.Range("I100").Value = 1146.2
$ws.Range("J100").Value = 1395
$ws.Range("K100").Value = 2292.4
$ws.Range("L100").Value = 2790
$ws.Range("M100").Value = -1751.4
$ws.Range("N100").Value = -3872

$ws.Range("H107").Value = 359
$ws.Range("I107").Value = 255.5
$ws.Range("K107").Value = 766.5
$ws.Range("M107").Value = 1153.5

$ws.Range("H113").Value = 742.7
$ws.Range("I113").Value = 742.7
$ws.Range("K113").Value = 2228.1
$ws.Range("M113").Value = -58.10000000000036

$ws.Range("H132").Value = 3471.0571
$ws.Range("I132").Value = 3269.7334
$ws.Range("K132").Value = 9809.200199999999
$ws.Range("M132").Value = -7279.200199999999

$ws.Range("H136").Value = 5351.2354
$ws.Range("J136").Value = 14633.333
$ws.Range("L136").Value = 43899.999
$ws.Range("N136").Value = -48999.999
